# CIV-11205 - split merge-field placeholder runs (adding proofErr markers
# Word inserts around mid-sentence grammar-checked fragments) and update the
# court/site placeholder text and the reasonAvailable conditional run.
$d = $word.ActiveDocument

# --- Change 1: "<<{dateFormat($nowUTC ,'d MMMM yyyy')}>>" run is split into
# three runs ("<", "<{", "dateFormat(...)}>>") with proofErr gramStart/gramEnd
# markers bracketing the middle run. ---
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("dateFormat(`$nowUTC", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) { throw "Change 1: anchor text not found" }
$para1 = $rng1.Paragraphs(1).Range
$xml1 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="0F2A3481" w14:textId="77777777" w:rsidR="00260D77" w:rsidRDefault="00260D77" w:rsidP="00C200CE"><w:pPr><w:rPr><w:rFonts w:ascii="GDS Transport Website Light" w:hAnsi="GDS Transport Website Light"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="GDS Transport Website Light" w:hAnsi="GDS Transport Website Light"/></w:rPr><w:t>&lt;</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="GDS Transport Website Light" w:hAnsi="GDS Transport Website Light"/></w:rPr><w:t>&lt;{</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="GDS Transport Website Light" w:hAnsi="GDS Transport Website Light"/></w:rPr><w:t>dateFormat($nowUTC ,‘d MMMM yyyy’)}&gt;&gt;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $para1.InsertXML($xml1)

# --- Change 2: the trailing "... on <<submittedOn>>>> at <<courtName>>." run is
# split so the placeholder becomes "<<siteName>> - <<address>> - <<postcode>>." ---
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("at <<courtName>>.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "Change 2: anchor text not found" }
$para2 = $rng2.Paragraphs(1).Range
$xml2 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="3C52EF22" w14:textId="77777777" w:rsidR="00260D77" w:rsidRDefault="00260D77" w:rsidP="00260D77"><w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>This order is made by &lt;&lt;judgeNameTitle&gt;&gt; on &lt;&lt;</w:t></w:r><w:r><w:t>submittedOn</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">&gt;&gt; at </w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/></w:rPr><w:t>&lt;&lt;siteName&gt;&gt; - &lt;&lt;address&gt;&gt; - &lt;&lt;postcode&gt;&gt;.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $para2.InsertXML($xml2)

# --- Change 3: "<<cs_{reasonAvailable=='Yes' }>> " run is split into three runs
# ("<<cs", "_{", "reasonAvailable=='Yes' }>> ") with proofErr markers around the
# middle run, matching the pattern used for Change 1. ---
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("reasonAvailable", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found3) { throw "Change 3: anchor text not found" }
$para3 = $rng3.Paragraphs(1).Range
$xml3 = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="75A31D65" w14:textId="77777777" w:rsidR="00260D77" w:rsidRDefault="00260D77" w:rsidP="00260D77"><w:pPr><w:rPr><w:color w:val="000000"/><w:sz w:val="27"/><w:szCs w:val="27"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="27"/><w:szCs w:val="27"/></w:rPr><w:t>&lt;&lt;cs</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="27"/><w:szCs w:val="27"/></w:rPr><w:t>_{</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="27"/><w:szCs w:val="27"/></w:rPr><w:t xml:space="preserve">reasonAvailable==’Yes’ }&gt;&gt; </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="000000"/><w:sz w:val="27"/><w:szCs w:val="27"/></w:rPr><w:t>REASONS:</w:t></w:r><w:r><w:rPr><w:color w:val="000000"/><w:sz w:val="27"/><w:szCs w:val="27"/></w:rPr><w:t xml:space="preserve"> &lt;&lt;es_ &gt;&gt;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $para3.InsertXML($xml3)

Write-Output "CIV-11205 edits applied"
